$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 7 = "Fish Restoration Program": add a Survey_link hyperlink to its Tidal Wetland monitoring page.
# Prime the cell with the Hyperlink cell style (+ vertical-center alignment, matching the rest of the
# table) before wiring up the link, so the final cell format matches the other linked cells in the sheet.
$rng7 = $ws.Range("B7")
$rng7.Style = "Hyperlink"
$rng7.VerticalAlignment = -4108  # xlCenter
$rng7.Borders.LineStyle = -4142  # xlLineStyleNone (re-affirm the (lack of) border)
$ws.Hyperlinks.Add($rng7, "https://iep.ca.gov/Science-Synthesis-Service/Monitoring-Programs/Tidal-Wetland", [Type]::Missing, [Type]::Missing, "https://iep.ca.gov/Science-Synthesis-Service/Monitoring-Programs/Tidal-Wetland")

# Row 4 = "Enhanced Delta Smelt Monitoring": add a Survey_link hyperlink to the FWS juvenile fish
# monitoring program page.
$rng4 = $ws.Range("B4")
$rng4.Style = "Hyperlink"
$rng4.VerticalAlignment = -4108  # xlCenter
$ws.Hyperlinks.Add($rng4, "https://www.fws.gov/lodi/juvenile_fish_monitoring_program/jfmp_index.htm", [Type]::Missing, [Type]::Missing, "https://www.fws.gov/lodi/juvenile_fish_monitoring_program/jfmp_index.htm")

# Match the exact formatting already used by the other hyperlink cells in column H (Hyperlink style +
# border + vertical-center, no fill).
$ws.Range("H14").Copy()
$rng4.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Move the active selection to H8, as recorded in the saved view state.
$ws.Range("H8").Select()
